# Auto-generated Excel COM-interop edit script
# Applies cell value updates to the cryptocurrency price/volume table
# (price refresh + three row-pair swaps: LEO/Filecoin, Toncoin/Dai, EnergySwap/Maker)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to store a value as literal text, bypassing Excels
# automatic number/date inference (needed for numeric-looking price strings
# such as "8.00" or "1.00" that must keep their original formatting),
# then restores the default "Normal" style so no stray number format remains.
function Set-TextValue {
    param($Address, $Text)
    $cell = $ws.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = '62.465.08'
$ws.Range("E2").Value = '  +9.61%  '

$ws.Range("D3").Value = '3.374.23'
$ws.Range("E3").Value = '  +4.15%  '

$ws.Range("E4").Value = '  +0.04%  '

Set-TextValue "D5" '415.41'
$ws.Range("E5").Value = '  +5.21%  '

Set-TextValue "D6" '117.50'
$ws.Range("E6").Value = '  +8.75%  '

$ws.Range("D7").Value = '3.371.28'
$ws.Range("E7").Value = '  +4.14%  '

Set-TextValue "D8" '0.575'
$ws.Range("E8").Value = '  -1.90%  '

$ws.Range("E9").Value = '  -0.11%  '

Set-TextValue "D11" '0.116'
$ws.Range("E11").Value = '  +18.30%  '

Set-TextValue "D12" '40.21'
$ws.Range("E12").Value = '  +2.58%  '

$ws.Range("E13").Value = '  -0.48%  '

$ws.Range("D14").Value = '3.901.96'
$ws.Range("E14").Value = '  +3.98%  '

Set-TextValue "D15" '8.35'
$ws.Range("E15").Value = '  -0.02%  '

Set-TextValue "D16" '19.41'
$ws.Range("E16").Value = '  +2.00%  '

$ws.Range("D17").Value = '3.398.77'
$ws.Range("E17").Value = '  +4.60%  '

$ws.Range("D18").Value = '62.253.53'
$ws.Range("E18").Value = '  +9.44%  '

$ws.Range("E19").Value = '  -1.63%  '

Set-TextValue "D20" '10.91'
$ws.Range("E20").Value = '  +0.45%  '

Set-TextValue "D21" '0.0000117'
$ws.Range("E21").Value = '  +7.20%  '

$ws.Range("E22").Value = '  +0.26%  '

Set-TextValue "D23" '12.62'
$ws.Range("E23").Value = '  -2.98%  '

Set-TextValue "D24" '296.75'
$ws.Range("E24").Value = '  +1.69%  '

Set-TextValue "D25" '74.67'
$ws.Range("E25").Value = '  +0.90%  '

Set-TextValue "D26" '3.13'
$ws.Range("E26").Value = '  -0.87%  '

Set-TextValue "D27" '29.59'
$ws.Range("E27").Value = '  +5.42%  '

Set-TextValue "D28" '8.00'

Set-TextValue "D29" '0.175'
$ws.Range("E29").Value = '  +3.76%  '

$ws.Range("B30").Value = 'Filecoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue "D30" '7.70'
$ws.Range("E30").Value = '  +0.64%  '

$ws.Range("B31").Value = 'LEO'
$ws.Range("C31").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue "D31" '4.27'
$ws.Range("E31").Value = '  -2.01%  '

Set-TextValue "D32" '43.46'
$ws.Range("E32").Value = '  +8.66%  '

$ws.Range("E33").Value = '  +4.37%  '

Set-TextValue "D34" '11.42'
$ws.Range("E34").Value = '  +2.00%  '

$ws.Range("B35").Value = 'Dai'
$ws.Range("C35").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue "D35" '0.999'
$ws.Range("E35").Value = '  +0.01%  '

$ws.Range("B36").Value = 'Toncoin'
$ws.Range("C36").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue "D36" '2.53'
$ws.Range("E36").Value = '  +18.78%  '

Set-TextValue "D37" '0.0487'
$ws.Range("E37").Value = '  -0.47%  '

Set-TextValue "D38" '52.31'
$ws.Range("E38").Value = '  +1.54%  '

Set-TextValue "D39" '3.10'
$ws.Range("E39").Value = '  +5.54%  '

Set-TextValue "D40" '1.00'
$ws.Range("E40").Value = '  +0.01%  '

Set-TextValue "D41" '3.41'
$ws.Range("E41").Value = '  -1.38%  '

Set-TextValue "D42" '133.61'
$ws.Range("E42").Value = '  -2.58%  '

$ws.Range("E43").Value = '  -1.48%  '

$ws.Range("E44").Value = '  +3.53%  '

Set-TextValue "D45" '1.90'
$ws.Range("E45").Value = '  +0.32%  '

Set-TextValue "D46" '3.88'
$ws.Range("E46").Value = '  -2.03%  '

$ws.Range("E47").Value = '  -2.91%  '

$ws.Range("E48").Value = '  -4.86%  '

$ws.Range("B49").Value = 'Maker'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D49").Value = '2.172.63'
$ws.Range("E49").Value = '  +0.84%  '

$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue "D50" '21.23'
$ws.Range("E50").Value = '  -4.36%  '

$ws.Range("D51").Value = '3.706.54'
$ws.Range("E51").Value = '  +3.88%  '
